# Upgraded for vs 2012
#
# 1) Refresh the cached "datetimeFigureOut" date placeholder text (the
#    slide master + every slide layout carries its own cached copy of
#    that field) from 08/12/2012 to 26/12/2012.
# 2) Append a new, blank slide (id 257) at the end of the deck.

$p = $ppt.ActivePresentation

$oldDate = "08/12/2012"
$newDate = "26/12/2012"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# --- Fix the date placeholder cached on the slide master itself ---
$master = $p.Slides.Item(1).Master
Update-DatePlaceholder $master.Shapes

# --- Fix the date placeholder cached on every slide layout ---
# PowerPoint exposes $master.CustomLayouts, but the only reliable way to
# reach each individual layout's real shapes in this object model is via
# a slide that actually uses that layout, so briefly create one slide per
# layout, patch its layout's placeholder, then remove the scratch slide
# again (this does not disturb the slide id counter).
$ppLayoutTitle          = 1
$ppLayoutText           = 2
$ppLayoutSectionHeader  = 33
$ppLayoutTwoContent     = 4
$ppLayoutComparison     = 5
$ppLayoutTitleOnly      = 6
$ppLayoutBlank          = 12
$ppLayoutContentCaption = 8
$ppLayoutPictureCaption = 9
$ppLayoutVerticalTitleAndText = 10
$ppLayoutVerticalText   = 28

$layoutsToFix = @(
    $ppLayoutTitle,
    $ppLayoutText,
    $ppLayoutSectionHeader,
    $ppLayoutTwoContent,
    $ppLayoutComparison,
    $ppLayoutTitleOnly,
    $ppLayoutBlank,
    $ppLayoutContentCaption,
    $ppLayoutPictureCaption,
    $ppLayoutVerticalTitleAndText,
    $ppLayoutVerticalText
)

foreach ($layoutId in $layoutsToFix) {
    $scratch = $p.Slides.Add($p.Slides.Count + 1, $layoutId)
    Update-DatePlaceholder $scratch.CustomLayout.Shapes
    $scratch.Delete()
}

# --- Append the new blank slide at the end of the presentation ---
$newSlide = $p.Slides.Add($p.Slides.Count + 1, $ppLayoutBlank)
